{"js": "// Update the date title and every arithmetic expression in the 20x5\n// practice table. Each table cell / the title paragraph holds exactly\n// one run, so we replace text at the paragraph level (not via\n// cell.body.insertText) to preserve the existing run/paragraph\n// formatting (rFonts, sz, jc, etc.) instead of resetting it.\n\nconst newTitle = \"2025-11-28 Friday\";\n\nconst newTableValues = [\n  [\"10+43=\", \"48+43=\", \"39-2=\", \"72-62=\", \"92-43=\"],\n  [\"49+8=\", \"20+8=\", \"12+10=\", \"99-10=\", \"30+60=\"],\n  [\"87-33=\", \"68-3=\", \"4+8=\", \"19+59=\", \"59+1=\"],\n  [\"38+10=\", \"82-25=\", \"55-41=\", \"81-18=\", \"80-22=\"],\n  [\"56-40=\", \"93-87=\", \"34-32=\", \"4+53=\", \"15+54=\"],\n  [\"97-23=\", \"61-53=\", \"54-48=\", \"53-10=\", \"4+55=\"],\n  [\"92-3=\", \"92-91=\", \"31+63=\", \"60+14=\", \"80+8=\"],\n  [\"68+5=\", \"91-86=\", \"36+35=\", \"45+32=\", \"33+3=\"],\n  [\"91-53=\", \"45+15=\", \"31+23=\", \"12+37=\", \"60-47=\"],\n  [\"39-27=\", \"90-46=\", \"99-32=\", \"67-67=\", \"16+14=\"],\n  [\"8+55=\", \"96-46=\", \"41-15=\", \"18+34=\", \"9+82=\"],\n  [\"51-21=\", \"2+32=\", \"3+50=\", \"47-42=\", \"46+52=\"],\n  [\"2+88=\", \"91-1=\", \"41+1=\", \"91-50=\", \"46-37=\"],\n  [\"7+72=\", \"48-37=\", \"68+5=\", \"38-0=\", \"17+57=\"],\n  [\"30-28=\", \"50+26=\", \"33+57=\", \"40+25=\", \"40+11=\"],\n  [\"65-45=\", \"11+8=\", \"94-60=\", \"65+0=\", \"45-41=\"],\n  [\"56+22=\", \"81-0=\", \"34-22=\", \"22+21=\", \"78-24=\"],\n  [\"52-7=\", \"96-24=\", \"29+42=\", \"23+38=\", \"63-14=\"],\n  [\"6-1=\", \"3+34=\", \"36-11=\", \"26-18=\", \"81-19=\"],\n  [\"12+42=\", \"18+67=\", \"90-23=\", \"74-35=\", \"87-81=\"],\n];\n\nconst body = context.document.body;\n\n// --- Title paragraph (the date line above the table) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.insertText(newTitle, Word.InsertLocation.replace);\n\n// --- The practice table (20 rows x 5 columns of expressions) ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nfor (let r = 0; r < rowCount; r++) {\n  const newRow = newTableValues[r];\n  for (let c = 0; c < newRow.length; c++) {\n    const cell = table.getCell(r, c);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n    const cellParagraph = cellParagraphs.items[0];\n    cellParagraph.insertText(newRow[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and every arithmetic expression in the 20x5\n# practice table. Each table cell / the title paragraph holds exactly\n# one run, so we set Range.Text directly (instead of deleting/recreating\n# the cell content), which preserves the existing run/paragraph\n# formatting (rFonts, sz, jc, etc.).\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph (the date line above the table) ---\n$d.Paragraphs.Item(1).Range.Text = \"2025-11-28 Friday\"\n\n# --- The practice table (20 rows x 5 columns of expressions) ---\n$newTableValues = @(\n    @(\"10+43=\", \"48+43=\", \"39-2=\", \"72-62=\", \"92-43=\"),\n    @(\"49+8=\", \"20+8=\", \"12+10=\", \"99-10=\", \"30+60=\"),\n    @(\"87-33=\", \"68-3=\", \"4+8=\", \"19+59=\", \"59+1=\"),\n    @(\"38+10=\", \"82-25=\", \"55-41=\", \"81-18=\", \"80-22=\"),\n    @(\"56-40=\", \"93-87=\", \"34-32=\", \"4+53=\", \"15+54=\"),\n    @(\"97-23=\", \"61-53=\", \"54-48=\", \"53-10=\", \"4+55=\"),\n    @(\"92-3=\", \"92-91=\", \"31+63=\", \"60+14=\", \"80+8=\"),\n    @(\"68+5=\", \"91-86=\", \"36+35=\", \"45+32=\", \"33+3=\"),\n    @(\"91-53=\", \"45+15=\", \"31+23=\", \"12+37=\", \"60-47=\"),\n    @(\"39-27=\", \"90-46=\", \"99-32=\", \"67-67=\", \"16+14=\"),\n    @(\"8+55=\", \"96-46=\", \"41-15=\", \"18+34=\", \"9+82=\"),\n    @(\"51-21=\", \"2+32=\", \"3+50=\", \"47-42=\", \"46+52=\"),\n    @(\"2+88=\", \"91-1=\", \"41+1=\", \"91-50=\", \"46-37=\"),\n    @(\"7+72=\", \"48-37=\", \"68+5=\", \"38-0=\", \"17+57=\"),\n    @(\"30-28=\", \"50+26=\", \"33+57=\", \"40+25=\", \"40+11=\"),\n    @(\"65-45=\", \"11+8=\", \"94-60=\", \"65+0=\", \"45-41=\"),\n    @(\"56+22=\", \"81-0=\", \"34-22=\", \"22+21=\", \"78-24=\"),\n    @(\"52-7=\", \"96-24=\", \"29+42=\", \"23+38=\", \"63-14=\"),\n    @(\"6-1=\", \"3+34=\", \"36-11=\", \"26-18=\", \"81-19=\"),\n    @(\"12+42=\", \"18+67=\", \"90-23=\", \"74-35=\", \"87-81=\"),\n)\n\n$table = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newTableValues.Count; $r++) {\n    $rowValues = $newTableValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $table.Rows.Item($r).Cells.Item($c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
